$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for columns L..S (row 1)
$ws.Range("L1").Value = "hzj-混合调节_20170516_152754_ASIC_EEG"
$ws.Range("M1").Value = "hzj-混合调节_20170518_134207_ASIC_EEG"
$ws.Range("N1").Value = "hzj-混合调节_20170519_135415_ASIC_EEG"
$ws.Range("O1").Value = "zyx-混合调节_20170516_111228_ASIC_EEG"
$ws.Range("P1").Value = "zyx-混合调节_20170517_110944_ASIC_EEG"
$ws.Range("Q1").Value = "zyx-混合调节_20170518_112337_ASIC_EEG"
$ws.Range("R1").Value = "zyx-混合调节_20170519_124954_ASIC_EEG"
$ws.Range("S1").Value = "zyx-混合调节_20170522_111557_ASIC_EEG"

# Row 2 values
$ws.Range("L2").Value = 0.99305555555555558
$ws.Range("M2").Value = 0.95189003436426112
$ws.Range("N2").Value = 0.92509363295880154
$ws.Range("O2").Value = 0.9358974358974359
$ws.Range("P2").Value = 0.93269230769230771
$ws.Range("Q2").Value = 0.93203883495145634
$ws.Range("R2").Value = 0.96451612903225803
$ws.Range("S2").Value = 0.93689320388349517

# Row 3 values
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.967741935483871
$ws.Range("N3").Value = 0.9926739926739927
$ws.Range("O3").Value = 0.97857142857142854
$ws.Range("P3").Value = 0.95530726256983245
$ws.Range("Q3").Value = 0.96907216494845361
$ws.Range("R3").Value = 0.94630872483221473
$ws.Range("S3").Value = 0.93506493506493504

# Update the selection to cover the full used range, matching the diff
$ws.Range("A1:S3").Select()
